# Insert a new data row right after the header row (row 1), shifting all
# existing data rows down by one, then populate the new row 2 with the
# Honeywell（霍尼韦尔）entry and keep the running index in column A in
# sequence (0, 1, 2, ...).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push everything currently in rows 2..32 down to rows 3..33.
$ws.Rows(2).Insert()

# The engine's row-insert carries over formatting from neighbouring rows;
# strip that so the new row matches the plain (unstyled) look of the other
# data rows.
$ws.Range("A2:S2").ClearFormats()

# Column A is a simple running index (0-based). After inserting the new
# row, bump every pre-existing data row's index up by one so it still
# reads 0, 1, 2, ... down the sheet.
for ($r = 33; $r -ge 3; $r--) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Honeywell（霍尼韦尔）"
$ws.Range("C2").Value = "江宁大学城"
$ws.Range("D2").Value = "软件"
$ws.Range("E2").Value = "测试开发工程师"
$ws.Range("F2").Value = "8:30-17:00"
$ws.Range("G2").Value = "11:30-13:00"
$ws.Range("H2").Value = "看情况，基本到点就走，除非自己判断deadline前无法做完，周末几乎无加班，周末加班调休"
$ws.Range("I2").Value = "工资全额的10%+年底补充年薪5%"
$ws.Range("J2").Value = "固定13薪，额外奖金1个月左右"
$ws.Range("K2").Value = "试用期6个月，不打折，转正目前没有答辩，也不要写申请"
$ws.Range("L2").Value = "后端测试Dell,前端MAC 工位人均3,4平，空间很大，能放个床"
$ws.Range("M2").Value = "12天年假"
$ws.Range("N2").Value = "不打卡"
$ws.Range("O2").Value = "工作节奏不快，员工福利基本上每个月都会发一些东西，其他也没什么了"
$ws.Range("Q2").Value = "2022-06-23 10:02:47"

# Restore the bold, bordered, centered style used by every other cell in
# column A (direct "Style =" assignment is a no-op in this engine, so set
# the underlying format properties instead).
$ws.Range("A2").Font.Bold = $true
$ws.Range("A2").HorizontalAlignment = -4108
$ws.Range("A2").VerticalAlignment = -4160
$ws.Range("A2").Borders.LineStyle = 1
